# Update workout values (camera class lower filter) and set active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Tricep press: 130 -> 150
$ws.Range("C4").Value = 150

# Leg extensions: 60 -> 70
$ws.Range("C9").Value = 70

# Biceps curl: 50 -> 60
$ws.Range("C15").Value = 60

# Update the selected/active cell shown in the sheet view
$ws.Activate()
$ws.Range("C16").Select()
